# Auto update: 2025-12-05 02:00:49
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Date column (A2:A4): "2025-12-03" -> "2025-12-05" ---
# Force text so Excel doesn't coerce the ISO-looking date string into a date serial.
$ws.Range("A2:A4").NumberFormat = "@"
$ws.Range("A2").Value = "2025-12-05"
$ws.Range("A3").Value = "2025-12-05"
$ws.Range("A4").Value = "2025-12-05"

# --- Contract label (B4): "Gold Dec 25" -> "Gold Feb 26" ---
$ws.Range("B4").Value = "Gold Feb 26"

# --- MACRO_SIGNAL column (O2:O4): bullish -> neutral ---
$ws.Range("O2").Value = "⚪ 중립 구간"
$ws.Range("O3").Value = "⚪ 중립 구간"
$ws.Range("O4").Value = "⚪ 중립 구간"

# --- Row 2 (GLD) numeric refresh ---
$ws.Range("D2").Value = 387.6
$ws.Range("E2").Value = 56.8
$ws.Range("F2").Value = 1.17
$ws.Range("H2").Value = 73
$ws.Range("I2").Value = 83
$ws.Range("J2").Value = 93
$ws.Range("K2").Value = 66.90000000000001
$ws.Range("N2").Value = 52.43913937059539

# --- Row 3 (NEM) numeric refresh ---
$ws.Range("D3").Value = 90.45
$ws.Range("E3").Value = 52.1
$ws.Range("F3").Value = -0.08
$ws.Range("G3").Value = 50
$ws.Range("H3").Value = 80
$ws.Range("I3").Value = 80
$ws.Range("J3").Value = 83
$ws.Range("K3").Value = 62.7
$ws.Range("N3").Value = 52.43913937059539

# --- Row 4 (GC=F / Gold Feb 26) numeric refresh ---
$ws.Range("D4").Value = 4243.3
$ws.Range("E4").Value = 71.8
$ws.Range("F4").Value = 4.48
$ws.Range("H4").Value = 43
$ws.Range("K4").Value = 52.9
$ws.Range("N4").Value = 52.43913937059539
